$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Standalone data correction on row 19 (DIMRA 20 F.C.TAB.) ---
$ws.Range("H19").Value = "1:0"
$ws.Range("L19").Value = 70
$ws.Range("N19").Value = "1:0"

# --- 2. Insert a new product row before row 46 (LIDOCAINE ...), shifting
#        rows 46:106 down to 47:107 ---
$ws.Rows("46:46").Insert()

# Pick up the formatting (styles/borders/fonts) for the new row from the
# row immediately below it (which now holds the old row-46 formatting).
$ws.Range("A47:N47").Copy()
$ws.Range("A46:N46").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Recreate the merged cell regions used by every product row.
$ws.Range("B46:G46").Merge()
$ws.Range("H46:K46").Merge()
$ws.Range("L46:M46").Merge()

# Match the row height used by the other (odd-positioned) product rows.
$ws.Rows("46:46").RowHeight = 25.5

# Populate the new row with the new product's data.
$ws.Range("A46").Value = 43
$ws.Range("B46").Value = "LACRITEARS EYE DROPS 15 ML"
$ws.Range("H46").Value = "1:0"
$ws.Range("L46").Value = 49
$ws.Range("N46").Value = "1:0"

# The "م" sequence number in column A is positional (row-3), not tied to
# the product, so it has to be re-stamped for every row the insert shifted
# (the Insert() call above dragged the old numbers down together with the
# rest of the row).
for ($r = 47; $r -le 105; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 3
}

# --- 3. Refresh the grand-total cell (now on row 106) to include the new
#        row and the row-19 correction ---
$ws.Range("K106").Value = 4348.6400000000003

Write-Host "edit applied"
